$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "2020" column (K) gets duplicated into a brand-new column L (same
# values + formatting), extending the table by one more year column.
$ws.Columns("K").Copy() | Out-Null
$ws.Columns("L").Insert(-4161) | Out-Null   # xlShiftToRight
$excel.CutCopyMode = 0

# Update the selected/active cell shown in the saved view.
$ws.Range("L10").Select() | Out-Null
